# Adiciona dois novos blocos de "dicionario de dados" (tabelas associativas
# filme_premio e filme_ator) ao final da planilha "Filme", replicando o
# padrao visual (cabecalho + linhas de dados com celula A mesclada + linha
# separadora) ja usado para as tabelas existentes (filme, diretor, ator,
# premio).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- largura das novas colunas A e B (antes ocultas/nao usadas nesta aba) ---
$ws.Columns("A").ColumnWidth = 23.83
$ws.Columns("B").ColumnWidth = 18

# ------------------------------------------------------------------
# Bloco 1: filme_premio (associativa) -> linhas 26-31
# ------------------------------------------------------------------

# linha separadora (mesmo visual da linha 6 / 13 / 19)
$ws.Range("A6:E6").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)

# linha de cabecalho (Tabela | Coluna | Tipo | Restricoes | Descricao)
$ws.Range("A1:E1").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Tabela"
$ws.Range("B27").Value = "Coluna"
$ws.Range("C27").Value = "Tipo"
$ws.Range("D27").Value = "Restrições"
$ws.Range("E27").Value = "Descrição"

# 4 linhas de dados (estilo igual ao das linhas 2-5 da tabela "filme")
$ws.Range("A2:E2").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$ws.Range("A4:E4").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)

$ws.Range("A28").Value = "filme_premio (associativa)"
$ws.Range("B28").Value = "id_filme_premio"
$ws.Range("C28").Value = "UUID"
$ws.Range("D28").Value = "PK, NOT NULL"
$ws.Range("E28").Value = "Identificador único do ator, gerado automaticamente."

$ws.Range("B29").Value = "id_filme"
$ws.Range("C29").Value = "UUID"
$ws.Range("D29").Value = "FK, NOT NULL"
$ws.Range("E29").Value = "Referência ao filme associado."

$ws.Range("B30").Value = "id_premio"
$ws.Range("C30").Value = "UUID"
$ws.Range("D30").Value = "FK, NOT NULL"
$ws.Range("E30").Value = "Referência ao prêmio associado."

$ws.Range("B31").Value = "ano_premio"
$ws.Range("C31").Value = "INTEGER"
$ws.Range("D31").Value = "NOT NULL"
$ws.Range("E31").Value = "Ano em que o prêmio foi recebido."

$ws.Range("A28:A31").Merge()

# ------------------------------------------------------------------
# Bloco 2: filme_ator (associativa) -> linhas 32-36
# ------------------------------------------------------------------

# linha separadora
$ws.Range("A6:E6").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)

# linha de cabecalho
$ws.Range("A1:E1").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Tabela"
$ws.Range("B33").Value = "Coluna"
$ws.Range("C33").Value = "Tipo"
$ws.Range("D33").Value = "Restrições"
$ws.Range("E33").Value = "Descrição"

# 3 linhas de dados
$ws.Range("A2:E2").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)

$ws.Range("A34").Value = "filme_ator (associativa)"
$ws.Range("B34").Value = "id_filme_ator"
$ws.Range("C34").Value = "UUID"
$ws.Range("D34").Value = "PK, NOT NULL"
$ws.Range("E34").Value = "Identificador único do ator, gerado automaticamente."

$ws.Range("B35").Value = "id_filme"
$ws.Range("C35").Value = "UUID"
$ws.Range("D35").Value = "FK, NOT NULL"
$ws.Range("E35").Value = "Referência ao filme associado."

$ws.Range("B36").Value = "id_ator"
$ws.Range("C36").Value = "UUID"
$ws.Range("D36").Value = "FK, NOT NULL"
$ws.Range("E36").Value = "Referência ao ator associado."

$ws.Range("A34:A36").Merge()

# ------------------------------------------------------------------
# seleciona a celula abaixo do novo conteudo, como no arquivo final
# ------------------------------------------------------------------
$ws.Range("B39").Select()

Write-Output "done"
